$wb = $excel.ActiveWorkbook

# Sheet "Мой_Сценарий" (3rd sheet) - remove the helper column D formulas,
# replace the C7 total with a legacy CSE array formula, and make this the
# active sheet/selection.
$wsScenario = $wb.Worksheets.Item(3)

$wsScenario.Range("D2:D5").ClearContents()
$wsScenario.Range("C7").FormulaArray = "=SUM(B2:B5*C2:C5)"

# Update the selection (drives tabSelected / activeTab / selection/activeCell)
$wsScenario.Range("F7").Select()
